$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8-26 down to 9-27
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = 10
$ws.Range("B8").Value = "Vega Modelo de Temuco"
$ws.Range("C8").Value = "La Araucanía"
$ws.Range("D8").Value = 44467
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 100112026
$ws.Range("G8").Value = "Haba"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 40
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 14000
$ws.Range("N8").Value = '$/saco 25 kilos'
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 560
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
